$d = $word.ActiveDocument

# 1. Fix the garbled OCR text in the first paragraph:
#    "24." -> "24," and the mis-scanned "8" -> "B" in the set-up list.
$d.Content.Find.Execute(
    "24. Benjamin carried out an experiment using four set-ups A, 8, C and D as shown",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "24, Benjamin carried out an experiment using four set-ups A, B, C and D as shown",
    2) | Out-Null

# 2. Clean up the stray "+ " before "heating? |" in the question paragraph
#    (this paragraph currently sits right after the first picture).
$d.Content.Find.Execute(
    "+ heating? |",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "heating? |",
    2) | Out-Null

# 3. Remove the first picture's paragraph entirely (picture + its paragraph
#    mark), which lets the "Inwhich set-up..." question paragraph take its
#    place right after the intro paragraph.
$firstPicPara = $d.Paragraphs.Item(2)
$d.Range($firstPicPara.Range.Start, $firstPicPara.Range.End).Delete()

# 4. Remove the second picture's paragraph entirely (it is now paragraph 3,
#    sitting right after the question paragraph).
$secondPicPara = $d.Paragraphs.Item(3)
$d.Range($secondPicPara.Range.Start, $secondPicPara.Range.End).Delete()

# 5. Fix the garbled OCR text in the final answer paragraph:
#    ". (3). c |" -> ". @).\u00a2c |"
$d.Content.Find.Execute(
    ". (3). c |",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". @)." + [char]0x00A2 + "c |",
    2) | Out-Null
